$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'29.212.38"
$ws.Cells.Item(2, 5).Value = '  +0.11%  '
$ws.Cells.Item(3, 4).Value = "'1.834.91"
$ws.Cells.Item(3, 5).Value = '  -0.43%  '
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = '  +0.20%  '
$ws.Cells.Item(5, 4).Value = "'240.41"
$ws.Cells.Item(5, 5).Value = '  -0.61%  '
$ws.Cells.Item(6, 4).Value = "'0.6652"
$ws.Cells.Item(6, 5).Value = '  -3.28%  '
$ws.Cells.Item(7, 4).Value = "'1.001"
$ws.Cells.Item(7, 5).Value = '  +0.11%  '
$ws.Cells.Item(8, 4).Value = "'0.07329"
$ws.Cells.Item(8, 5).Value = '  -2.07%  '
$ws.Cells.Item(9, 4).Value = "'0.2908"
$ws.Cells.Item(9, 5).Value = '  -3.61%  '
$ws.Cells.Item(10, 4).Value = "'22.58"
$ws.Cells.Item(10, 5).Value = '  -2.64%  '
$ws.Cells.Item(11, 4).Value = "'0.07679"
$ws.Cells.Item(11, 5).Value = '  +0.22%  '
$ws.Cells.Item(12, 4).Value = "'1.829.35"
$ws.Cells.Item(12, 5).Value = '  -1.12%  '
$ws.Cells.Item(13, 4).Value = "'4.954"
$ws.Cells.Item(13, 5).Value = '  -2.26%  '
$ws.Cells.Item(14, 4).Value = "'0.6643"
$ws.Cells.Item(14, 5).Value = '  -2.94%  '
$ws.Cells.Item(15, 4).Value = "'82.95"
$ws.Cells.Item(15, 5).Value = '  -5.36%  '
$ws.Cells.Item(16, 4).Value = "'6.085"
$ws.Cells.Item(16, 5).Value = '  -1.67%  '
$ws.Cells.Item(17, 4).Value = "'29.195.09"
$ws.Cells.Item(17, 5).Value = '  +0.12%  '
$ws.Cells.Item(18, 4).Value = "'0.000008248"
$ws.Cells.Item(18, 5).Value = '  +0.71%  '
$ws.Cells.Item(19, 4).Value = "'225.13"
$ws.Cells.Item(19, 5).Value = '  -2.05%  '
$ws.Cells.Item(20, 4).Value = "'12.43"
$ws.Cells.Item(20, 5).Value = '  -1.12%  '
$ws.Cells.Item(21, 5).Value = '  +0.17%  '
$ws.Cells.Item(22, 4).Value = "'7.109"
$ws.Cells.Item(22, 5).Value = '  -4.05%  '
$ws.Cells.Item(23, 4).Value = "'1.002"
$ws.Cells.Item(23, 5).Value = '  +0.31%  '
$ws.Cells.Item(24, 4).Value = "'160.39"
$ws.Cells.Item(24, 5).Value = '  +0.55%  '
$ws.Cells.Item(25, 4).Value = "'8.615"
$ws.Cells.Item(25, 5).Value = '  -2.02%  '
$ws.Cells.Item(26, 4).Value = "'0.1388"
$ws.Cells.Item(26, 5).Value = '  -4.63%  '
$ws.Cells.Item(27, 4).Value = "'17.90"
$ws.Cells.Item(27, 5).Value = '  -1.21%  '
$ws.Cells.Item(28, 4).Value = "'1.510"
$ws.Cells.Item(28, 5).Value = '  -0.05%  '
$ws.Cells.Item(29, 4).Value = "'4.101"
$ws.Cells.Item(29, 5).Value = '  -4.46%  '
$ws.Cells.Item(30, 4).Value = "'4.026"
$ws.Cells.Item(30, 5).Value = '  -3.03%  '
$ws.Cells.Item(31, 4).Value = "'1.193"
$ws.Cells.Item(31, 5).Value = '  -0.23%  '
$ws.Cells.Item(32, 4).Value = "'0.05274"
$ws.Cells.Item(32, 5).Value = '  +0.55%  '
$ws.Cells.Item(33, 4).Value = "'1.862"
$ws.Cells.Item(33, 5).Value = '  +0.31%  '
$ws.Cells.Item(34, 4).Value = "'0.7466"
$ws.Cells.Item(34, 5).Value = '  -1.97%  '
$ws.Cells.Item(35, 4).Value = "'1.126"
$ws.Cells.Item(35, 5).Value = '  -0.98%  '
$ws.Cells.Item(36, 4).Value = "'2.681"
$ws.Cells.Item(36, 5).Value = '  -0.04%  '
$ws.Cells.Item(37, 4).Value = "'1.313.91"
$ws.Cells.Item(37, 5).Value = '  +0.54%  '
$ws.Cells.Item(38, 4).Value = "'0.01789"
$ws.Cells.Item(38, 5).Value = '  -2.55%  '
$ws.Cells.Item(39, 4).Value = "'2.716"
$ws.Cells.Item(39, 5).Value = '  -0.39%  '
$ws.Cells.Item(40, 4).Value = "'0.9178"
$ws.Cells.Item(40, 5).Value = '  -1.87%  '
$ws.Cells.Item(41, 2).Value = 'XinFinNetwork'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Cells.Item(41, 4).Value = "'0.08608"
$ws.Cells.Item(41, 5).Value = '  +18.93%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).Value = "'5.938"
$ws.Cells.Item(42, 5).Value = '  -0.51%  '
$ws.Cells.Item(43, 5).Value = '  +0.29%  '
$ws.Cells.Item(44, 4).Value = "'101.66"
$ws.Cells.Item(44, 5).Value = '  -3.38%  '
$ws.Cells.Item(45, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(45, 4).Value = "'0.00000000125"
$ws.Cells.Item(45, 5).Value = '  +1.98%  '
$ws.Cells.Item(46, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(46, 4).Value = "'1.979.50"
$ws.Cells.Item(46, 5).Value = '  -0.28%  '
$ws.Cells.Item(47, 2).Value = 'Mantle'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(47, 4).Value = "'0.5167"
$ws.Cells.Item(47, 5).Value = '  -0.54%  '
$ws.Cells.Item(48, 2).Value = 'RenderToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(48, 4).Value = "'1.760"
$ws.Cells.Item(48, 5).Value = '  -0.99%  '
$ws.Cells.Item(49, 4).Value = "'63.51"
$ws.Cells.Item(49, 5).Value = '  -2.69%  '
$ws.Cells.Item(50, 4).Value = "'9.070"
$ws.Cells.Item(50, 5).Value = '  -4.55%  '
$ws.Cells.Item(51, 4).Value = "'0.05935"
$ws.Cells.Item(51, 5).Value = '  -0.30%  '
